$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Tweak lesson #37 name (drop "example" from the Rest server lesson title)
$ws.Range("C41").Value2 = "Rest server – glassfish, final project"

# 2) Fill in the date + youtube link for lesson #38 (row 42) which were blank before
$ws.Range("E42").Value2 = 44242
$ws.Hyperlinks.Add($ws.Range("F42"), "https://youtu.be/tVDrwZ32lOk", "", "", "https://youtu.be/tVDrwZ32lOk")
$ws.Range("F42").Value2 = "https://youtu.be/tVDrwZ32lOk "
# Hyperlinks.Add reformats the cell with the blue/underlined hyperlink
# style; restore the plain formatting used by the rest of column F here.
$ws.Range("F40").Copy()
$ws.Range("F42").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) New section header "F. Spring" + lesson #39 "Spring intro" (row 43)
$ws.Range("A43").Value2 = "F. Spring"
$ws.Range("C43").Value2 = "Spring intro"
$ws.Range("D43").Value2 = 2
$ws.Range("E43").Value2 = 44243

# 4) Lesson #40 "Spring intro" continuation (row 44)
$ws.Range("C44").Value2 = "Spring intro"
$ws.Range("D44").Value2 = 2
$ws.Range("E44").Value2 = 44246

# 5) Row-height autofit touch-ups to mirror Excel's automatic recalculation
$ws.Rows.Item(41).RowHeight = 14.9
$ws.Rows.Item(43).RowHeight = 14.9
$ws.Rows.Item(44).RowHeight = 13.85

# 6) Selection state left by the author after the edit
$ws.Range("C46").Select()
